# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Mon Jan 29 20:59:43 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.138.31"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "'2.302.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'310.62"
$ws.Range("D6").Value = "'101.25"
$ws.Range("E6").Value = "  +6.28%  "
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +7.56%  "
$ws.Range("D10").Value = "'36.02"
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("E11").Value = "  +3.93%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "'7.17"
$ws.Range("E13").Value = "  +7.53%  "
$ws.Range("D14").Value = "'2.657.24"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "'2.296.30"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "'43.070.70"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("D19").Value = "'12.55"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "'68.48"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "'240.97"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("E24").Value = "  +2.68%  "
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "'24.75"
$ws.Range("E27").Value = "  +4.78%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  +10.13%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.18"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").Value = "'168.14"
$ws.Range("E31").Value = "  +4.94%  "
$ws.Range("D32").Value = "'5.32"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'17.75"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "'4.23"
$ws.Range("E41").Value = "  +5.73%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "'0.0290"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").Value = "'1.977.36"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "'19.17"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").Value = "'9.86"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'2.97"
$ws.Range("E48").Value = "  +18.81%  "
$ws.Range("D49").Value = "'55.53"
$ws.Range("E49").Value = "  +4.43%  "
$ws.Range("D50").Value = "'2.525.55"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "  +2.55%  "
